$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''30.517.97'
$ws.Range("E2").Value = '  +0.37%  '
$ws.Range("D3").Value = '''1.852.05'
$ws.Range("E3").Value = '  -0.40%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '''233.52'
$ws.Range("E5").Value = '  -0.64%  '
$ws.Range("E6").Value = '  +0.00%  '
$ws.Range("D7").Value = '''0.4702'
$ws.Range("E7").Value = '  -0.91%  '
$ws.Range("D8").Value = '''0.2747'
$ws.Range("E8").Value = '  -0.07%  '
$ws.Range("D9").Value = '''0.06339'
$ws.Range("E9").Value = '  -1.60%  '
$ws.Range("D10").Value = '''17.60'
$ws.Range("E10").Value = '  +7.64%  '
$ws.Range("D11").Value = '''1.836.93'
$ws.Range("E11").Value = '  -1.36%  '
$ws.Range("D12").Value = '''0.07419'
$ws.Range("E12").Value = '  -0.26%  '
$ws.Range("D13").Value = '''5.031'
$ws.Range("E13").Value = '  +0.46%  '
$ws.Range("D14").Value = '''84.59'
$ws.Range("E14").Value = '  -1.32%  '
$ws.Range("D15").Value = '''0.6266'
$ws.Range("E15").Value = '  -1.40%  '
$ws.Range("D16").Value = '''30.495.29'
$ws.Range("E16").Value = '  +0.45%  '
$ws.Range("D17").Value = '''243.11'
$ws.Range("E17").Value = '  +5.12%  '
$ws.Range("D18").Value = '''0.9999'
$ws.Range("E18").Value = '  +0.04%  '
$ws.Range("D19").Value = '''12.67'
$ws.Range("E19").Value = '  -1.11%  '
$ws.Range("D20").Value = '''0.000007341'
$ws.Range("E20").Value = '  -1.11%  '
$ws.Range("D21").Value = '''1.001'
$ws.Range("E21").Value = '  +0.02%  '
$ws.Range("D22").Value = '''4.955'
$ws.Range("E22").Value = '  -1.26%  '
$ws.Range("D23").Value = '''5.981'
$ws.Range("E23").Value = '  -0.70%  '
$ws.Range("D24").Value = '''9.261'
$ws.Range("E24").Value = '  +0.24%  '
$ws.Range("D25").Value = '''162.77'
$ws.Range("E25").Value = '  -1.95%  '
$ws.Range("D26").Value = '''18.06'
$ws.Range("E26").Value = '  +0.35%  '
$ws.Range("D27").Value = '''1.885'
$ws.Range("E27").Value = '  -0.71%  '
$ws.Range("E28").Value = '  -2.63%  '
$ws.Range("D29").Value = '''1.372'
$ws.Range("E29").Value = '  -1.61%  '
$ws.Range("D30").Value = '''4.041'
$ws.Range("E30").Value = '  -2.81%  '
$ws.Range("D31").Value = '''3.857'
$ws.Range("E31").Value = '  -2.09%  '
$ws.Range("D32").Value = '''0.04889'
$ws.Range("E32").Value = '  -0.64%  '
$ws.Range("D33").Value = '''1.139'
$ws.Range("E33").Value = '  -1.43%  '
$ws.Range("D34").Value = '''0.7044'
$ws.Range("E34").Value = '  -3.36%  '
$ws.Range("E35").Value = '  +0.32%  '
$ws.Range("D36").Value = '''0.01904'
$ws.Range("E36").Value = '  +0.31%  '
$ws.Range("D37").Value = '''2.682'
$ws.Range("E37").Value = '  +1.29%  '
$ws.Range("D38").Value = '''0.8718'
$ws.Range("E38").Value = '  -4.75%  '
$ws.Range("D39").Value = '''1.977'
$ws.Range("E39").Value = '  +0.18%  '
$ws.Range("D40").Value = '''105.33'
$ws.Range("E40").Value = '  -0.65%  '
$ws.Range("D41").Value = '''0.9999'
$ws.Range("E41").Value = '  +0.01%  '
$ws.Range("D42").Value = '''0.4075'
$ws.Range("E42").Value = '  -1.12%  '
$ws.Range("D43").Value = '''5.494'
$ws.Range("E43").Value = '  -1.66%  '
$ws.Range("D44").Value = '''7.240'
$ws.Range("E44").Value = '  +1.30%  '
$ws.Range("D45").Value = '''62.76'
$ws.Range("E45").Value = '  +2.51%  '
$ws.Range("E46").Value = '  -1.13%  '
$ws.Range("B47").Value = 'Elrond'
$ws.Range("C47").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D47").Value = '''33.36'
$ws.Range("E47").Value = '  -0.50%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '''8.530'
$ws.Range("E48").Value = '  -2.24%  '
$ws.Range("B49").Value = 'Cronos'
$ws.Range("C49").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D49").Value = '''0.05532'
$ws.Range("E49").Value = '  -1.01%  '
$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").Value = '''1.370'
$ws.Range("E50").Value = '  -2.81%  '
$ws.Range("D51").Value = '''0.3689'
$ws.Range("E51").Value = '  -0.76%  '
